$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 new rows before row 99; this shifts the existing rows 99:206
# down to 103:210 (the previous rows 203:206 become the new rows 207:210).
$ws.Rows.Item(99).Resize(4).Insert()

# Shared/constant values for every data row in this block.
$mercadoId = 5
$mercado   = "Macroferia Regional de Talca"
$region    = "Maule"
$codreg    = 7
$catId     = 100112027
$categoria = "Mel$([char]0x00F3)n"
$unidad    = "`$/unidad"
$origen    = "Regi$([char]0x00F3)n del Maule"
$kgUnidad  = 1
$clasif    = "Hortaliza"

function Set-DataRow {
    param($row, $fecha, $variedad, $calidad, $volumen, $precio)

    $ws.Range("A$row").Value = $mercadoId
    $ws.Range("B$row").Value = $mercado
    $ws.Range("C$row").Value = $region
    $ws.Range("D$row").Value = $fecha
    $ws.Range("E$row").Value = $codreg
    $ws.Range("F$row").Value = $catId
    $ws.Range("G$row").Value = $categoria
    $ws.Range("H$row").Value = $variedad
    $ws.Range("I$row").Value = $calidad
    $ws.Range("J$row").Value = $volumen
    $ws.Range("K$row").Value = $precio
    $ws.Range("L$row").Value = $precio
    $ws.Range("M$row").Value = $precio
    $ws.Range("N$row").Value = $unidad
    $ws.Range("O$row").Value = $origen
    $ws.Range("P$row").Value = $precio
    $ws.Range("Q$row").Value = $kgUnidad
    $ws.Range("R$row").Value = $clasif
}

# New rows 99-102: a new reporting date (44557) for Calame$([char]0x00F1)o / Tuna, Extra & Primera.
Set-DataRow -row 99  -fecha 44557 -variedad "Calame$([char]0x00F1)o" -calidad "Extra"   -volumen 3000 -precio 1000
Set-DataRow -row 100 -fecha 44557 -variedad "Calame$([char]0x00F1)o" -calidad "Primera" -volumen 6000 -precio 800
Set-DataRow -row 101 -fecha 44557 -variedad "Tuna"                    -calidad "Extra"   -volumen 3000 -precio 1000
Set-DataRow -row 102 -fecha 44557 -variedad "Tuna"                    -calidad "Primera" -volumen 6000 -precio 800
